# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting a newer data snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 8306
$wsExpo.Range("F5").Value = 35757
$wsExpo.Range("F12").Value = 817
$wsExpo.Range("F15").Value = 449
$wsExpo.Range("F20").Value = 429
$wsExpo.Range("F24").Value = 2401
$wsExpo.Range("F30").Value = 678
$wsExpo.Range("F31").Value = 678

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 356

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 8307
$wsAll.Range("F7").Value = 35757
$wsAll.Range("F15").Value = 356
$wsAll.Range("F18").Value = 817
$wsAll.Range("F21").Value = 449
$wsAll.Range("F31").Value = 429
$wsAll.Range("F35").Value = 2401
$wsAll.Range("F42").Value = 678
$wsAll.Range("F43").Value = 678
